# Actualización automática 2025-06-13 15:20:08
#
# Registers a sale of 798.31 for client "SANCHEZ CORREA MARCO EDUARDO"
# (group "240X80 PORCELANATO", month "junio") and a sale of 543.38 for
# client "MAD&DECO S.A." (group "FREGADEROS DE COCINA", month "junio"),
# and refreshes the dependent totals/percentages across all three sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" : per-client sales broken down by group
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# MAD&DECO S.A. (row 29) -> FREGADEROS DE COCINA (col E)
$wsGrupo.Range("E29").Value = 543.38

# SANCHEZ CORREA MARCO EDUARDO (row 44) -> 240X80 PORCELANATO (col D)
$wsGrupo.Range("D44").Value = 798.3099999999999

# Totals row: "n de 53" counters for columns D and E go up by one each
# now that a previously-zero cell in each column became non-zero.
$wsGrupo.Range("D55").Value = "7 de 53"
$wsGrupo.Range("E55").Value = "2 de 53"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL" : per-client sales broken down by month
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# MAD&DECO S.A. (row 29) -> junio (col F)
$wsMensual.Range("F29").Value = 543.38

# SANCHEZ CORREA MARCO EDUARDO (row 44) -> junio (col F)
$wsMensual.Range("F44").Value = 798.3099999999999

# Totals row for junio
$wsMensual.Range("F55").Value = 34957.74

# Column D ("abril") narrows slightly (from stored width 14 to 13).
# COM's ColumnWidth is offset from the raw OOXML "width" attribute by the
# fixed ~0.8333 default-font padding, so feed it the inverse value.
$wsMensual.Columns.Item(4).ColumnWidth = 12.166666666666666

# ---------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL" : completion % per group
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# 240X80 PORCELANATO (row 3)
$wsCumpl.Range("D3").Value = 12870.1
$wsCumpl.Range("E3").Value = 857.8999999999996
$wsCumpl.Range("F3").Value = 0.9375072843822844

# FREGADEROS DE COCINA (row 4)
$wsCumpl.Range("D4").Value = 672.12
$wsCumpl.Range("E4").Value = -26.12
$wsCumpl.Range("F4").Value = 1.040433436532508

# TOTAL row (19)
$wsCumpl.Range("D19").Value = 35243.99
$wsCumpl.Range("E19").Value = 55719.339
$wsCumpl.Range("F19").Value = 0.3874527283406701
